$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: fill column A (ID) for all 30 new rows, in row order
$ws.Cells.Item(1172, 1).Value = "U40_01"
$ws.Cells.Item(1173, 1).Value = "U40_02"
$ws.Cells.Item(1174, 1).Value = "U40_03"
$ws.Cells.Item(1175, 1).Value = "U40_04"
$ws.Cells.Item(1176, 1).Value = "U40_05"
$ws.Cells.Item(1177, 1).Value = "U40_06"
$ws.Cells.Item(1178, 1).Value = "U40_07"
$ws.Cells.Item(1179, 1).Value = "U40_08"
$ws.Cells.Item(1180, 1).Value = "U40_09"
$ws.Cells.Item(1181, 1).Value = "U40_10"
$ws.Cells.Item(1182, 1).Value = "U40_11"
$ws.Cells.Item(1183, 1).Value = "U40_12"
$ws.Cells.Item(1184, 1).Value = "U40_13"
$ws.Cells.Item(1185, 1).Value = "U40_14"
$ws.Cells.Item(1186, 1).Value = "U40_15"
$ws.Cells.Item(1187, 1).Value = "U40_16"
$ws.Cells.Item(1188, 1).Value = "U40_17"
$ws.Cells.Item(1189, 1).Value = "U40_18"
$ws.Cells.Item(1190, 1).Value = "U40_19"
$ws.Cells.Item(1191, 1).Value = "U40_20"
$ws.Cells.Item(1192, 1).Value = "U40_21"
$ws.Cells.Item(1193, 1).Value = "U40_22"
$ws.Cells.Item(1194, 1).Value = "U40_23"
$ws.Cells.Item(1195, 1).Value = "U40_24"
$ws.Cells.Item(1196, 1).Value = "U40_25"
$ws.Cells.Item(1197, 1).Value = "U40_26"
$ws.Cells.Item(1198, 1).Value = "U40_27"
$ws.Cells.Item(1199, 1).Value = "U40_28"
$ws.Cells.Item(1200, 1).Value = "U40_29"
$ws.Cells.Item(1201, 1).Value = "U40_30"

# Phase 2: fill columns B..G per row in the specific order matching the source data entry
# Row 1172
$ws.Cells.Item(1172, 2).Value = 40
$ws.Cells.Item(1172, 3).Value = "Cánh"
$ws.Cells.Item(1172, 4).Value = "Wing"
$ws.Cells.Item(1172, 6).Value = "chicken wing"
$ws.Cells.Item(1172, 5).Value = "KFC is known for its (nổi tiếng với) chicken wings"
$ws.Cells.Item(1172, 7).Value = "N"

# Row 1173
$ws.Cells.Item(1173, 2).Value = 40
$ws.Cells.Item(1173, 3).Value = "Chạm trán, gặp phải"
$ws.Cells.Item(1173, 4).Value = "Encounter"
$ws.Cells.Item(1173, 5).Value = "The detective(thám tử) encountered a hard case"
$ws.Cells.Item(1173, 6).Value = "encounter somebody or something"
$ws.Cells.Item(1173, 7).Value = "V"

# Row 1174
$ws.Cells.Item(1174, 2).Value = 40
$ws.Cells.Item(1174, 3).Value = "Sư tử"
$ws.Cells.Item(1174, 4).Value = "Lion"
$ws.Cells.Item(1174, 5).Value = "The lion cub follows its mother to learn hunting"
$ws.Cells.Item(1174, 6).Value = "lion cub / sư tử con"
$ws.Cells.Item(1174, 7).Value = "N"

# Row 1175
$ws.Cells.Item(1175, 2).Value = 40
$ws.Cells.Item(1175, 3).Value = "Vương quốc"
$ws.Cells.Item(1175, 4).Value = "Kingdom"
$ws.Cells.Item(1175, 5).Value = "The kingdom of Portugal(Bồ Đào Nha) does not exist anymore"
$ws.Cells.Item(1175, 6).Value = "The kingdom of "
$ws.Cells.Item(1175, 7).Value = "N"

# Row 1176
$ws.Cells.Item(1176, 2).Value = 40
$ws.Cells.Item(1176, 3).Value = "Khổng lồ"
$ws.Cells.Item(1176, 4).Value = "Giant"
$ws.Cells.Item(1176, 5).Value = "The aquarium(thủy cung) has a giant squid"
$ws.Cells.Item(1176, 6).Value = "a giant squid / một con mực khổng lồ"
$ws.Cells.Item(1176, 7).Value = "Adj"

# Row 1177
$ws.Cells.Item(1177, 2).Value = 40
$ws.Cells.Item(1177, 3).Value = "Hoang dã, hoang dại"
$ws.Cells.Item(1177, 4).Value = "Wild"
$ws.Cells.Item(1177, 5).Value = "Wild animals are hard to tame (thuần hóa)"
$ws.Cells.Item(1177, 6).Value = "Wild animals / động vật hoang dã"
$ws.Cells.Item(1177, 7).Value = "Adj"

# Row 1178
$ws.Cells.Item(1178, 2).Value = 40
$ws.Cells.Item(1178, 3).Value = "Kẻ thù"
$ws.Cells.Item(1178, 4).Value = "Enemy"
$ws.Cells.Item(1178, 5).Value = "She is my sworn enemy"
$ws.Cells.Item(1178, 6).Value = "Somebody's sworn enemy / kẻ thù không đội trời chung"
$ws.Cells.Item(1178, 7).Value = "N"

# Row 1179
$ws.Cells.Item(1179, 2).Value = 40
$ws.Cells.Item(1179, 3).Value = "Săn"
$ws.Cells.Item(1179, 4).Value = "Hunt"
$ws.Cells.Item(1179, 6).Value = "hunt for something"
$ws.Cells.Item(1179, 5).Value = "People in that tribe(bộ tộc) hunt for deer (hươu)"
$ws.Cells.Item(1179, 7).Value = "V"

# Row 1180
$ws.Cells.Item(1180, 2).Value = 40
$ws.Cells.Item(1180, 3).Value = "Lồng, chuồng"
$ws.Cells.Item(1180, 4).Value = "Cage"
$ws.Cells.Item(1180, 5).Value = "He doesn't want to keep the dog in a cage."
$ws.Cells.Item(1180, 6).Value = "keep something in a cage"
$ws.Cells.Item(1180, 7).Value = "N"

# Row 1181
$ws.Cells.Item(1181, 2).Value = 40
$ws.Cells.Item(1181, 3).Value = "Rộng rãi"
$ws.Cells.Item(1181, 4).Value = "Spacious"
$ws.Cells.Item(1181, 5).Value = "The spacious room allowed us to dance"
$ws.Cells.Item(1181, 6).Value = "a spacious room"
$ws.Cells.Item(1181, 7).Value = "Adj"

# Row 1182
$ws.Cells.Item(1182, 2).Value = 40
$ws.Cells.Item(1182, 3).Value = "Con dê"
$ws.Cells.Item(1182, 4).Value = "goat"
$ws.Cells.Item(1182, 5).Value = "There are no mountain goats here"
$ws.Cells.Item(1182, 6).Value = "a mountain goat / dê núi"
$ws.Cells.Item(1182, 7).Value = "N"

# Row 1183
$ws.Cells.Item(1183, 2).Value = 40
$ws.Cells.Item(1183, 3).Value = "Gãi, cào"
$ws.Cells.Item(1183, 4).Value = "Scratch"
$ws.Cells.Item(1183, 5).Value = "The cat scratched my face"
$ws.Cells.Item(1183, 6).Value = "scratch something"
$ws.Cells.Item(1183, 7).Value = "V"

# Row 1184
$ws.Cells.Item(1184, 2).Value = 40
$ws.Cells.Item(1184, 3).Value = "Chủng, giống"
$ws.Cells.Item(1184, 4).Value = "Breed"
$ws.Cells.Item(1184, 5).Value = "The farm has many breeds of sheep"
$ws.Cells.Item(1184, 6).Value = "A breed of sheep / giống cừu"
$ws.Cells.Item(1184, 7).Value = "N"

# Row 1185
$ws.Cells.Item(1185, 2).Value = 40
$ws.Cells.Item(1185, 3).Value = "Lông"
$ws.Cells.Item(1185, 4).Value = "Fur"
$ws.Cells.Item(1185, 5).Value = "I wear a fur coat because the weather is cold"
$ws.Cells.Item(1185, 6).Value = "a fur coat / áo khoác lông thú"
$ws.Cells.Item(1185, 7).Value = "N"

# Row 1186
$ws.Cells.Item(1186, 2).Value = 40
$ws.Cells.Item(1186, 3).Value = "Thô cứng"
$ws.Cells.Item(1186, 4).Value = "Rough"
$ws.Cells.Item(1186, 5).Value = "My father is used (đã quen) to driving on rough roads"
$ws.Cells.Item(1186, 6).Value = "rough roads / những con đường gồ ghề"
$ws.Cells.Item(1186, 7).Value = "Adj"

# Row 1187
$ws.Cells.Item(1187, 2).Value = 40
$ws.Cells.Item(1187, 3).Value = "Thông minh"
$ws.Cells.Item(1187, 4).Value = "Intelligent"
$ws.Cells.Item(1187, 5).Value = "He student are highly intelligent"
$ws.Cells.Item(1187, 6).Value = "Highly intelligent / cực kỳ thông minh"
$ws.Cells.Item(1187, 7).Value = "Adj"

# Row 1188
$ws.Cells.Item(1188, 2).Value = 40
$ws.Cells.Item(1188, 3).Value = "Môi trường sống"
$ws.Cells.Item(1188, 4).Value = "Habitat"
$ws.Cells.Item(1188, 5).Value = "The jungle (rừng già) is the monkey's natural habitat"
$ws.Cells.Item(1188, 6).Value = "natural habitat / môi trường sống tự  nhiên"
$ws.Cells.Item(1188, 7).Value = "N"

# Row 1189
$ws.Cells.Item(1189, 2).Value = 40
$ws.Cells.Item(1189, 3).Value = "Cỏ"
$ws.Cells.Item(1189, 4).Value = "Grass"
$ws.Cells.Item(1189, 5).Value = "Laying on (nằm trên) a field of grass is nice"
$ws.Cells.Item(1189, 6).Value = "a field of grass"
$ws.Cells.Item(1189, 7).Value = "N"

# Row 1190
$ws.Cells.Item(1190, 2).Value = 40
$ws.Cells.Item(1190, 3).Value = "Bẩn"
$ws.Cells.Item(1190, 4).Value = "Dirty"
$ws.Cells.Item(1190, 5).Value = "I have dirty hands after painting"
$ws.Cells.Item(1190, 6).Value = "dirty hands / tay bẩn"
$ws.Cells.Item(1190, 7).Value = "Adj"

# Row 1191
$ws.Cells.Item(1191, 2).Value = 40
$ws.Cells.Item(1191, 3).Value = "Dọn dẹp"
$ws.Cells.Item(1191, 4).Value = "Clean up"
$ws.Cells.Item(1191, 5).Value = "Clean this place up before he arrives"
$ws.Cells.Item(1191, 6).Value = "clean something up"
$ws.Cells.Item(1191, 7).Value = "V"

# Row 1192
$ws.Cells.Item(1192, 2).Value = 40
$ws.Cells.Item(1192, 3).Value = "Nhận nuôi"
$ws.Cells.Item(1192, 4).Value = "Adopt"
$ws.Cells.Item(1192, 5).Value = "They are going to adopt a child"
$ws.Cells.Item(1192, 6).Value = "adopt somebody / nhận nuôi một ai đó"
$ws.Cells.Item(1192, 7).Value = "V"

# Row 1193
$ws.Cells.Item(1193, 2).Value = 40
$ws.Cells.Item(1193, 3).Value = "Thú nuôi"
$ws.Cells.Item(1193, 4).Value = "Pet"
$ws.Cells.Item(1193, 5).Value = "Our family pet is a parrot (con vẹt)"
$ws.Cells.Item(1193, 6).Value = "a family pet / thú nuôi trong gia đình"
$ws.Cells.Item(1193, 7).Value = "N"

# Row 1194
$ws.Cells.Item(1194, 2).Value = 40
$ws.Cells.Item(1194, 3).Value = "Sinh học"
$ws.Cells.Item(1194, 4).Value = "Biology"
$ws.Cells.Item(1194, 5).Value = "I graduated with a degree in biology"
$ws.Cells.Item(1194, 6).Value = "a degree in biology / một tấm bằng ngành sinh học"
$ws.Cells.Item(1194, 7).Value = "N"

# Row 1195
$ws.Cells.Item(1195, 2).Value = 40
$ws.Cells.Item(1195, 3).Value = "Tận hưởng, yêu thích"
$ws.Cells.Item(1195, 4).Value = "Enjoy"
$ws.Cells.Item(1195, 5).Value = "Hope you enjoy the movie"
$ws.Cells.Item(1195, 6).Value = "enjoy something"
$ws.Cells.Item(1195, 7).Value = "V"

# Row 1196
$ws.Cells.Item(1196, 2).Value = 40
$ws.Cells.Item(1196, 3).Value = "Đua"
$ws.Cells.Item(1196, 4).Value = "Race"
$ws.Cells.Item(1196, 5).Value = "Five horses will race against each other in the Grand Finale (vòng chung kết)"
$ws.Cells.Item(1196, 6).Value = "race against / đua với"
$ws.Cells.Item(1196, 7).Value = "V"

# Row 1197
$ws.Cells.Item(1197, 2).Value = 40
$ws.Cells.Item(1197, 3).Value = "Yên lặng, trật tự"
$ws.Cells.Item(1197, 4).Value = "Silent"
$ws.Cells.Item(1197, 5).Value = "You must keep silent as everyone is asleep (đang ngủ)"
$ws.Cells.Item(1197, 6).Value = "keep silent"
$ws.Cells.Item(1197, 7).Value = "Adj"

# Row 1198
$ws.Cells.Item(1198, 2).Value = 40
$ws.Cells.Item(1198, 3).Value = "Điều chỉnh"
$ws.Cells.Item(1198, 4).Value = "Adjust"
$ws.Cells.Item(1198, 5).Value = "Can you adjust the seat belt for me?"
$ws.Cells.Item(1198, 6).Value = "adjust something"
$ws.Cells.Item(1198, 7).Value = "V"

# Row 1199
$ws.Cells.Item(1199, 2).Value = 40
$ws.Cells.Item(1199, 3).Value = "Đói"
$ws.Cells.Item(1199, 4).Value = "Hungry"
$ws.Cells.Item(1199, 5).Value = "He gives the hungry crowd enough bread and fish"
$ws.Cells.Item(1199, 6).Value = "a hungry crowd / một đám người đang đói"
$ws.Cells.Item(1199, 7).Value = "Adj"

# Row 1200
$ws.Cells.Item(1200, 2).Value = 40
$ws.Cells.Item(1200, 3).Value = "Cho ăn"
$ws.Cells.Item(1200, 4).Value = "Feed"
$ws.Cells.Item(1200, 5).Value = "My mother is feeding the baby"
$ws.Cells.Item(1200, 6).Value = "feed somebody"
$ws.Cells.Item(1200, 7).Value = "V"

# Row 1201
$ws.Cells.Item(1201, 2).Value = 40
$ws.Cells.Item(1201, 3).Value = "Hết sạch"
$ws.Cells.Item(1201, 4).Value = "Run out of"
$ws.Cells.Item(1201, 5).Value = "We ran out of milk"
$ws.Cells.Item(1201, 6).Value = "run out of something"
$ws.Cells.Item(1201, 7).Value = "V"

# Update the view state to match final selection
$ws.Application.ActiveWindow.ScrollRow = 1185
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("G1201").Select()